# Auto-generated Excel COM-interop script to apply the diff changes
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1310.6666  # H15 was 1376.375
$ws.Cells.Item(15, 9).Value = 1310.6666  # I15 was 1376.375
$ws.Cells.Item(15, 11).Value = 3931.9998  # K15 was 4129.125
$ws.Cells.Item(15, 13).Value = -3762.9998  # M15 was -3960.125
$ws.Cells.Item(19, 8).Value = 2010.88  # H19 was 2151.4546
$ws.Cells.Item(19, 10).Value = 1711.125  # J19 was 2149.8
$ws.Cells.Item(19, 12).Value = 1711.125  # L19 was 2149.8
$ws.Cells.Item(19, 14).Value = -2061.125  # N19 was -2499.8
$ws.Cells.Item(112, 8).Value = 4082.1614  # H112 was 4162.2666
$ws.Cells.Item(112, 10).Value = 4249.8965  # J112 was 4341.7144
$ws.Cells.Item(112, 12).Value = 12749.6895  # L112 was 13025.1432
$ws.Cells.Item(112, 14).Value = -14965.6895  # N112 was -15241.1432
$ws.Cells.Item(135, 8).Value = 3853.85  # H135 was 3761.4546
$ws.Cells.Item(135, 9).Value = 945.7857  # I135 was 1182.25
$ws.Cells.Item(135, 11).Value = 8512.0713  # K135 was 10640.25
$ws.Cells.Item(135, 13).Value = -5977.0713  # M135 was -8105.25
$ws.Cells.Item(141, 8).Value = 6044.5586  # H141 was 6372.4375
$ws.Cells.Item(141, 9).Value = 5400.6895  # I141 was 5741.593
$ws.Cells.Item(141, 11).Value = 16202.0685  # K141 was 17224.779
$ws.Cells.Item(141, 13).Value = -11022.0685  # M141 was -12044.779

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14183.151  # H32 was 14621.375
$ws.Cells.Item(32, 9).Value = 13388.29  # I32 was 13829.233
$ws.Cells.Item(32, 11).Value = 13388.29  # K32 was 13829.233
$ws.Cells.Item(32, 13).Value = -13101.29  # M32 was -13542.233
$ws.Cells.Item(43, 9).Value = 0  # I43 was 29999
$ws.Cells.Item(43, 10).Value = 29999.5  # J43 was 30000
$ws.Cells.Item(43, 11).Value = 0  # K43 was 29999
$ws.Cells.Item(43, 12).Value = 29999.5  # L43 was 30000
$ws.Cells.Item(43, 13).ClearContents()  # M43 removed (was -29686)
$ws.Cells.Item(43, 14).Value = -30625.5  # N43 was -30626
$ws.Cells.Item(61, 8).Value = 4429.533  # H61 was 3874.64
$ws.Cells.Item(61, 9).Value = 3536.9167  # I61 was 3052.4375
$ws.Cells.Item(61, 10).Value = 8000  # J61 was 5336.3335
$ws.Cells.Item(61, 11).Value = 3536.9167  # K61 was 3052.4375
$ws.Cells.Item(61, 12).Value = 8000  # L61 was 5336.3335
$ws.Cells.Item(61, 13).Value = -3324.9167  # M61 was -2840.4375
$ws.Cells.Item(61, 14).Value = -8424  # N61 was -5760.3335
$ws.Cells.Item(102, 8).Value = 2299.8333  # H102 was 2418.9092
$ws.Cells.Item(102, 9).Value = 2299.8333  # I102 was 2418.9092
$ws.Cells.Item(102, 11).Value = 2299.8333  # K102 was 2418.9092
$ws.Cells.Item(102, 13).Value = -677.8332999999998  # M102 was -796.9092000000001
$ws.Cells.Item(122, 8).Value = 5381350  # H122 was 5957596.5
$ws.Cells.Item(122, 9).Value = 6670674  # I122 was 7250596
$ws.Cells.Item(122, 10).Value = 9166.333000000001  # J122 was 9799.6
$ws.Cells.Item(122, 11).Value = 20012022  # K122 was 21751788
$ws.Cells.Item(122, 12).Value = 27498.999  # L122 was 29398.8
$ws.Cells.Item(122, 13).Value = -20009572  # M122 was -21749338
$ws.Cells.Item(122, 14).Value = -32398.999  # N122 was -34298.8
$ws.Cells.Item(136, 8).Value = 4429.533  # H136 was 3874.64
$ws.Cells.Item(136, 9).Value = 3536.9167  # I136 was 3052.4375
$ws.Cells.Item(136, 10).Value = 8000  # J136 was 5336.3335
$ws.Cells.Item(136, 11).Value = 10610.7501  # K136 was 9157.3125
$ws.Cells.Item(136, 12).Value = 24000  # L136 was 16009.0005
$ws.Cells.Item(136, 13).Value = -8060.750100000001  # M136 was -6607.3125
$ws.Cells.Item(136, 14).Value = -29100  # N136 was -21109.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(11, 8).Value = 75.583336  # H11 was 46.363636
$ws.Cells.Item(11, 9).Value = 15.6  # I11 was 50.9
$ws.Cells.Item(11, 10).Value = 375.5  # J11 was 1
$ws.Cells.Item(11, 11).Value = 15.6  # K11 was 50.9
$ws.Cells.Item(11, 12).Value = 375.5  # L11 was 1
$ws.Cells.Item(11, 13).Value = 124.4  # M11 was 89.09999999999999
$ws.Cells.Item(11, 14).Value = -655.5  # N11 was -281
$ws.Cells.Item(134, 8).Value = 2448.9092  # H134 was 2740.16
$ws.Cells.Item(134, 9).Value = 2043.75  # I134 was 2148.611
$ws.Cells.Item(134, 10).Value = 6500.5  # J134 was 4261.2856
$ws.Cells.Item(134, 11).Value = 6131.25  # K134 was 6445.833
$ws.Cells.Item(134, 12).Value = 19501.5  # L134 was 12783.8568
$ws.Cells.Item(134, 13).Value = -3596.25  # M134 was -3910.833
$ws.Cells.Item(134, 14).Value = -24571.5  # N134 was -17853.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1453.5  # H16 was 1453.4
$ws.Cells.Item(16, 9).Value = 1013.4  # I16 was 1013.2
$ws.Cells.Item(16, 11).Value = 1013.4  # K16 was 1013.2
$ws.Cells.Item(16, 13).Value = -726.4  # M16 was -726.2
$ws.Cells.Item(86, 8).Value = 10070.5  # H86 was 10311.667
$ws.Cells.Item(86, 9).Value = 9569.4  # I86 was 9769.4
$ws.Cells.Item(86, 10).Value = 10571.6  # J86 was 10989.5
$ws.Cells.Item(86, 11).Value = 9569.4  # K86 was 9769.4
$ws.Cells.Item(86, 12).Value = 10571.6  # L86 was 10989.5
$ws.Cells.Item(86, 13).Value = -8446.4  # M86 was -8646.4
$ws.Cells.Item(86, 14).Value = -12817.6  # N86 was -13235.5
$ws.Cells.Item(89, 8).Value = 10070.5  # H89 was 10311.667
$ws.Cells.Item(89, 9).Value = 9569.4  # I89 was 9769.4
$ws.Cells.Item(89, 10).Value = 10571.6  # J89 was 10989.5
$ws.Cells.Item(89, 11).Value = 47847  # K89 was 48847
$ws.Cells.Item(89, 12).Value = 52858  # L89 was 54947.5
$ws.Cells.Item(89, 13).Value = -42231  # M89 was -43231
$ws.Cells.Item(89, 14).Value = -64090  # N89 was -66179.5
$ws.Cells.Item(99, 8).Value = 6147.2666  # H99 was 5950.5
$ws.Cells.Item(99, 9).Value = 2713.625  # I99 was 2745.3333
$ws.Cells.Item(99, 11).Value = 2713.625  # K99 was 2745.3333
$ws.Cells.Item(99, 13).Value = -1215.625  # M99 was -1247.3333
$ws.Cells.Item(105, 8).Value = 1942  # H105 was 1999.25
$ws.Cells.Item(105, 10).Value = 1919.4  # J105 was 2000
$ws.Cells.Item(105, 12).Value = 1919.4  # L105 was 2000
$ws.Cells.Item(105, 14).Value = -5413.4  # N105 was -5494
$ws.Cells.Item(113, 8).Value = 1453.5  # H113 was 1453.4
$ws.Cells.Item(113, 9).Value = 1013.4  # I113 was 1013.2
$ws.Cells.Item(113, 11).Value = 1013.4  # K113 was 1013.2
$ws.Cells.Item(113, 13).Value = 1156.6  # M113 was 1156.8
$ws.Cells.Item(122, 8).Value = 2123.476  # H122 was 2123.5715
$ws.Cells.Item(122, 9).Value = 1480.5  # I122 was 1480.6333
$ws.Cells.Item(122, 11).Value = 4441.5  # K122 was 4441.8999
$ws.Cells.Item(122, 13).Value = -1991.5  # M122 was -1991.8999
$ws.Cells.Item(126, 8).Value = 6147.2666  # H126 was 5950.5
$ws.Cells.Item(126, 9).Value = 2713.625  # I126 was 2745.3333
$ws.Cells.Item(126, 11).Value = 8140.875  # K126 was 8235.999899999999
$ws.Cells.Item(126, 13).Value = -5670.875  # M126 was -5765.999899999999
$ws.Cells.Item(132, 8).Value = 49391372  # H132 was 55564870
$ws.Cells.Item(132, 9).Value = 63497100  # I132 was 74079384
$ws.Cells.Item(132, 11).Value = 190491300  # K132 was 222238152
$ws.Cells.Item(132, 13).Value = -190488770  # M132 was -222235622
$ws.Cells.Item(134, 8).Value = 4230  # H134 was 3020.6316
$ws.Cells.Item(134, 9).Value = 3922  # I134 was 2500.1333
$ws.Cells.Item(134, 10).Value = 5000  # J134 was 4972.5
$ws.Cells.Item(134, 11).Value = 11766  # K134 was 7500.3999
$ws.Cells.Item(134, 12).Value = 15000  # L134 was 14917.5
$ws.Cells.Item(134, 13).Value = -9231  # M134 was -4965.3999
$ws.Cells.Item(134, 14).Value = -20070  # N134 was -19987.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 835  # H18 was 566.6667
$ws.Cells.Item(18, 9).Value = 695  # I18 was 362.5
$ws.Cells.Item(18, 11).Value = 2085  # K18 was 1087.5
$ws.Cells.Item(18, 13).Value = -1916  # M18 was -918.5
$ws.Cells.Item(109, 8).Value = 7502.8887  # H109 was 8190.75
$ws.Cells.Item(109, 9).Value = 4705.2  # I109 was 5381.5
$ws.Cells.Item(109, 11).Value = 14115.6  # K109 was 16144.5
$ws.Cells.Item(109, 13).Value = -13075.6  # M109 was -15104.5
$ws.Cells.Item(114, 8).Value = 2127.5715  # H114 was 1688.421
$ws.Cells.Item(114, 10).Value = 2981.6667  # J114 was 2080.6428
$ws.Cells.Item(114, 12).Value = 8945.000100000001  # L114 was 6241.928400000001
$ws.Cells.Item(114, 14).Value = -15453.0001  # N114 was -12749.9284
$ws.Cells.Item(121, 8).Value = 67482.664  # H121 was 824.93335
$ws.Cells.Item(121, 9).Value = 250111.5  # I121 was 149
$ws.Cells.Item(121, 10).Value = 1072.1818  # J121 was 993.9167
$ws.Cells.Item(121, 11).Value = 750334.5  # K121 was 447
$ws.Cells.Item(121, 12).Value = 3216.5454  # L121 was 2981.7501
$ws.Cells.Item(121, 13).Value = -749024.5  # M121 was 863
$ws.Cells.Item(121, 14).Value = -5836.5454  # N121 was -5601.7501
$ws.Cells.Item(129, 8).Value = 1517.6923  # H129 was 1551.4073
$ws.Cells.Item(129, 10).Value = 3157.1428  # J129 was 3066
$ws.Cells.Item(129, 12).Value = 9471.428400000001  # L129 was 9198
$ws.Cells.Item(129, 14).Value = -19471.4284  # N129 was -19198

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 740.0714  # H97 was 785.46155
$ws.Cells.Item(97, 9).Value = 186.3  # I97 was 190.33333
$ws.Cells.Item(97, 11).Value = 186.3  # K97 was 190.33333
$ws.Cells.Item(97, 13).Value = 309.7  # M97 was 305.66667
$ws.Cells.Item(102, 8).Value = 15631376  # H102 was 16135604
$ws.Cells.Item(102, 9).Value = 20839086  # I102 was 21745098
$ws.Cells.Item(102, 10).Value = 8246  # J102 was 8308.5
$ws.Cells.Item(102, 11).Value = 20839086  # K102 was 21745098
$ws.Cells.Item(102, 12).Value = 8246  # L102 was 8308.5
$ws.Cells.Item(102, 13).Value = -20837464  # M102 was -21743476
$ws.Cells.Item(102, 14).Value = -11490  # N102 was -11552.5
$ws.Cells.Item(113, 8).Value = 3774096.5  # H113 was 18867924
$ws.Cells.Item(113, 9).Value = 4717606  # I113 was 18867924
$ws.Cells.Item(113, 10).Value = 60  # J113 was 0
$ws.Cells.Item(113, 11).Value = 4717606  # K113 was 18867924
$ws.Cells.Item(113, 12).Value = 60  # L113 was 0
$ws.Cells.Item(113, 13).Value = -4715436  # M113 was -18865754
$ws.Cells.Item(113, 14).Value = -4400  # N113 was N/A
$ws.Cells.Item(126, 8).Value = 3868.5417  # H126 was 3785.76
$ws.Cells.Item(126, 9).Value = 1930.2727  # I126 was 1919.3334
$ws.Cells.Item(126, 11).Value = 5790.8181  # K126 was 5758.0002
$ws.Cells.Item(126, 13).Value = -3320.8181  # M126 was -3288.0002
$ws.Cells.Item(132, 8).Value = 3736.5518  # H132 was 3898.7036
$ws.Cells.Item(132, 9).Value = 3278.4783  # I132 was 3443.3333
$ws.Cells.Item(132, 11).Value = 9835.4349  # K132 was 10329.9999
$ws.Cells.Item(132, 13).Value = -7305.4349  # M132 was -7799.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3827.3125  # H7 was 3703.4138
$ws.Cells.Item(7, 9).Value = 3206.2642  # I7 was 3099.0667
$ws.Cells.Item(7, 11).Value = 3206.2642  # K7 was 3099.0667
$ws.Cells.Item(7, 13).Value = -3094.2642  # M7 was -2987.0667
$ws.Cells.Item(16, 8).Value = 3287.125  # H16 was 2934.111
$ws.Cells.Item(16, 9).Value = 1216  # I16 was 1058
$ws.Cells.Item(16, 11).Value = 1216  # K16 was 1058
$ws.Cells.Item(16, 13).Value = -1046  # M16 was -888
$ws.Cells.Item(20, 8).Value = 7336666.5  # H20 was 6290700
$ws.Cells.Item(20, 9).Value = 20000  # I20 was 17449.5
$ws.Cells.Item(20, 11).Value = 20000  # K20 was 17449.5
$ws.Cells.Item(20, 13).Value = -19774  # M20 was -17223.5
$ws.Cells.Item(46, 8).Value = 3096.4546  # H46 was 3052.2778
$ws.Cells.Item(46, 9).Value = 2170.2307  # I46 was 2129.4285
$ws.Cells.Item(46, 10).Value = 3698.5  # J46 was 3639.5454
$ws.Cells.Item(46, 11).Value = 2170.2307  # K46 was 2129.4285
$ws.Cells.Item(46, 12).Value = 3698.5  # L46 was 3639.5454
$ws.Cells.Item(46, 13).Value = -1982.2307  # M46 was -1941.4285
$ws.Cells.Item(46, 14).Value = -4074.5  # N46 was -4015.5454
$ws.Cells.Item(122, 8).Value = 6859.3335  # H122 was 7188.5713
$ws.Cells.Item(122, 9).Value = 3333.5386  # I122 was 3530.5454
$ws.Cells.Item(122, 11).Value = 10000.6158  # K122 was 10591.6362
$ws.Cells.Item(122, 13).Value = -7550.6158  # M122 was -8141.636200000001
$ws.Cells.Item(126, 8).Value = 3827.3125  # H126 was 3703.4138
$ws.Cells.Item(126, 9).Value = 3206.2642  # I126 was 3099.0667
$ws.Cells.Item(126, 11).Value = 9618.792600000001  # K126 was 9297.2001
$ws.Cells.Item(126, 13).Value = -7148.792600000001  # M126 was -6827.2001
$ws.Cells.Item(127, 8).Value = 250068750  # H127 was 333408320
$ws.Cells.Item(127, 10).Value = 91666.336  # J127 was 112500
$ws.Cells.Item(127, 12).Value = 91666.336  # L127 was 112500
$ws.Cells.Item(127, 14).Value = -101586.336  # N127 was -122420
$ws.Cells.Item(128, 8).Value = 37000  # H128 was 0
$ws.Cells.Item(128, 10).Value = 37000  # J128 was 0
$ws.Cells.Item(128, 12).Value = 37000  # L128 was 0
$ws.Cells.Item(128, 14).Value = -46960  # N128 was N/A
$ws.Cells.Item(132, 8).Value = 2896.1848  # H132 was 2767.0312
$ws.Cells.Item(132, 9).Value = 2554.6875  # I132 was 2537.9397
$ws.Cells.Item(132, 10).Value = 5172.8335  # J132 was 4229.6924
$ws.Cells.Item(132, 11).Value = 7664.0625  # K132 was 7613.8191
$ws.Cells.Item(132, 12).Value = 15518.5005  # L132 was 12689.0772
$ws.Cells.Item(132, 13).Value = -5134.0625  # M132 was -5083.8191
$ws.Cells.Item(132, 14).Value = -20578.5005  # N132 was -17749.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(68, 8).Value = 36000  # H68 was 0
$ws.Cells.Item(68, 10).Value = 36000  # J68 was 0
$ws.Cells.Item(68, 12).Value = 36000  # L68 was 0
$ws.Cells.Item(68, 14).Value = -37622  # N68 was N/A
$ws.Cells.Item(71, 8).Value = 36000  # H71 was 0
$ws.Cells.Item(71, 10).Value = 36000  # J71 was 0
$ws.Cells.Item(71, 12).Value = 108000  # L71 was 0
$ws.Cells.Item(71, 14).Value = -116112  # N71 was N/A
$ws.Cells.Item(81, 8).Value = 5836.077  # H81 was 5874.231
$ws.Cells.Item(81, 10).Value = 6864.1665  # J81 was 6946.8335
$ws.Cells.Item(81, 12).Value = 13728.333  # L81 was 13893.667
$ws.Cells.Item(81, 14).Value = -15850.333  # N81 was -16015.667
$ws.Cells.Item(84, 8).Value = 5836.077  # H84 was 5874.231
$ws.Cells.Item(84, 10).Value = 6864.1665  # J84 was 6946.8335
$ws.Cells.Item(84, 12).Value = 68641.66500000001  # L84 was 69468.33499999999
$ws.Cells.Item(84, 14).Value = -79249.66500000001  # N84 was -80076.33499999999
$ws.Cells.Item(100, 8).Value = 1208.3125  # H100 was 1024.6666
$ws.Cells.Item(100, 9).Value = 718  # I100 was 721.63635
$ws.Cells.Item(100, 10).Value = 2287  # J100 was 1858
$ws.Cells.Item(100, 11).Value = 1436  # K100 was 1443.2727
$ws.Cells.Item(100, 12).Value = 4574  # L100 was 3716
$ws.Cells.Item(100, 13).Value = -895  # M100 was -902.2727
$ws.Cells.Item(100, 14).Value = -5656  # N100 was -4798
$ws.Cells.Item(107, 8).Value = 1403.6666  # H107 was 1433.375
$ws.Cells.Item(107, 9).Value = 1366.6  # I107 was 1416.75
$ws.Cells.Item(107, 11).Value = 4099.799999999999  # K107 was 4250.25
$ws.Cells.Item(107, 13).Value = -2179.799999999999  # M107 was -2330.25
$ws.Cells.Item(132, 8).Value = 1514.8909  # H132 was 1419.5333
$ws.Cells.Item(132, 9).Value = 2351.0667  # I132 was 1911.579
$ws.Cells.Item(132, 10).Value = 1201.325  # J132 was 1191.5122
$ws.Cells.Item(132, 11).Value = 7053.2001  # K132 was 5734.737
$ws.Cells.Item(132, 12).Value = 3603.975  # L132 was 3574.536599999999
$ws.Cells.Item(132, 13).Value = -4523.2001  # M132 was -3204.737
$ws.Cells.Item(132, 14).Value = -8663.975  # N132 was -8634.536599999999
